# upto data provider -3
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Add Customer")
$ws3 = $wb.Worksheets.Item("Sheet3")

# -----------------------------------------------------------------
# Sheet1: just a selection change (A1:B4 selected instead of B1:whole column)
# -----------------------------------------------------------------
$ws1.Range("A1:B4").Select()

# -----------------------------------------------------------------
# Add Customer sheet: bold the header row, switch the numeric columns
# back to General number format, and move the selection.
# -----------------------------------------------------------------
$ws2.Range("A1:D1").Font.Bold = $true
$ws2.Range("B2:C3").NumberFormat = "General"
$ws2.Range("E6").Select()

# -----------------------------------------------------------------
# Sheet3: fill in the login data-table (user / pass / expected url)
# with hyperlinks pointing at the "wrong password" error page.
# -----------------------------------------------------------------
$ws3.Range("A1").Value = "user"
$ws3.Range("B1").Value = "pass"
$ws3.Range("C1").Value = "expected url"
$ws3.Range("A1:C1").Font.Bold = $true

$ws3.Range("A2").Value = "admin"
$ws3.Range("B2").Value = "admin"
$ws3.Range("C2").Value = "http://stock.scriptinglogic.net/dashboard.php"

$ws3.Range("A3").Value = "excelX-1"
$ws3.Range("B3").Value = "wewe"
$ws3.Range("C3").Value = "http://stock.scriptinglogic.net/index.php?msg=Wrong%20Username%20or%20Password&type=error"

$ws3.Range("A4").Value = "excelX-2"
$ws3.Range("B4").Value = 1234
$ws3.Range("C4").Value = "http://stock.scriptinglogic.net/index.php?msg=Wrong%20Username%20or%20Password&type=error"

$ws3.Range("A5").Value = "excelX-3"
$ws3.Range("B5").Value = "xyz"
$ws3.Range("C5").Value = "http://stock.scriptinglogic.net/index.php?msg=Wrong%20Username%20or%20Password&type=error"

$ws3.Hyperlinks.Add($ws3.Range("C3"), "http://stock.scriptinglogic.net/index.php?msg=Wrong%20Username%20or%20Password&type=error")
$ws3.Hyperlinks.Add($ws3.Range("C4:C5"), "http://stock.scriptinglogic.net/index.php?msg=Wrong%20Username%20or%20Password&type=error", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "http://stock.scriptinglogic.net/index.php?msg=Wrong%20Username%20or%20Password&type=error")

$ws3.Columns.Item(3).ColumnWidth = 13.57
$ws3.Range("A1:C5").Select()
$ws3.Select()
$excel.ActiveWindow.Zoom = 235

# -----------------------------------------------------------------
# New "Sheet2" tab at the end: same user/pass table but checking the
# expected page title instead of a redirect URL (no hyperlinks yet).
# -----------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws4.Name = "Sheet2"

$ws4.Range("A1").Value = "user"
$ws4.Range("B1").Value = "pass"
$ws4.Range("C1").Value = "expected Page title"
$ws4.Range("A1:C1").Font.Bold = $true

$ws4.Range("A2").Value = "admin"
$ws4.Range("B2").Value = "admin"

$ws4.Range("A3").Value = "excelX-1"
$ws4.Range("B3").Value = "wewe"
$ws4.Range("C3").Value = ""

$ws4.Range("A4").Value = "excelX-2"
$ws4.Range("B4").Value = 1234
$ws4.Range("C4").Value = ""

$ws4.Range("A5").Value = "excelX-3"
$ws4.Range("B5").Value = "xyz"
$ws4.Range("C5").Value = ""

$ws4.Columns.Item(3).ColumnWidth = 17.3
$ws4.Range("C2:C5").Select()
$ws4.Activate()
$excel.ActiveWindow.Zoom = 205

Write-Host "edit complete"
